# Finalized tileset builder for 3x3 minimal tileset
#
# Three small "coordinate label" textboxes on slide 1 are updated:
#   - TextBox 246 (id 247) @ (5181600,1828800): "1,3 d" -> "1,3 " + "c", width 389850 -> 378630 EMU
#   - TextBox 249 (id 250) @ (5181600,3200400): "1,3 d" -> "1,3 " + "c", width 389850 -> 378630 EMU
#   - TextBox 309 (id 310) @ (7391400,4343400): "0,0 1" -> "0,0 a",      width 386644 -> 385042 EMU
#
# Shape.Width/Height/Left/Top are backed by single-precision (float32)
# storage, and EMU<-points conversion on save truncates rather than
# rounds. Add a tiny (< 1 EMU in point-units) epsilon before assigning
# so the value truncates back to the exact target EMU instead of
# landing 1 EMU short.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

# Splitting a run via TextRange.Characters(start,len).Text = "..." keeps
# each new run's rPr equal to the original run's rPr (lang/sz/dirty/
# smtClean all carried over), which is exactly the behavior the target
# markup needs for the "1,3 d" -> "1,3 "/"c" split.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TextBox 246 (id 247) ---------------------------------------------
$sh1 = $s.Shapes.Item("TextBox 246")
$sh1.Width = EmuToPt 378630
$sh1.TextFrame.TextRange.Characters(5, 1).Text = "c"

# --- TextBox 249 (id 250) ---------------------------------------------
$sh2 = $s.Shapes.Item("TextBox 249")
$sh2.Width = EmuToPt 378630
$sh2.TextFrame.TextRange.Characters(5, 1).Text = "c"

# --- TextBox 309 (id 310) ---------------------------------------------
$sh3 = $s.Shapes.Item("TextBox 309")
$sh3.Width = EmuToPt 385042
$sh3.TextFrame.TextRange.Text = "0,0 a"
